$wb = $excel.ActiveWorkbook

# --- Sheet "A18": append two new TxHash rows (A4, A5) ---
$wsA18 = $wb.Worksheets.Item("A18")
$wsA18.Range("A4").Value = "9252ACCE0D4F9620438287F16DF7E767A0DCC04F9E6DF5375D650EB5847C7259"
$wsA18.Range("A5").Value = "0F61C8C9A724D55250FCEFBB14435E2E9CCEFFC78CB9ED6D16C44860FA238FCC"
$wsA18.Range("B15").Select()

# --- Sheet "A19": append three new TxHash rows (A5, A6, A7) ---
$wsA19 = $wb.Worksheets.Item("A19")
$wsA19.Range("A5").Value = "74A9F95201E1F221BC6471DC3E074EF8081AD292B81A7BB352C0492D4D40B8F5"
$wsA19.Range("A6").Value = "EA5F39440473CB214BD042B87648AC5DA959BCF55B7AE1D0D9CD9FAFF3932812"
$wsA19.Range("A7").Value = "7647AE1D09B75F8245F17E14378E77FB15287F6907F03878968C8B3343269FF5"

# --- Sheet "A20": append three new TxHash rows (A5, A6, A7) ---
$wsA20 = $wb.Worksheets.Item("A20")
$wsA20.Range("A5").Value = "029B00403406C02819427D594BC715C52884CECD672B1EB542B3F2D686FC7766"
$wsA20.Range("A6").Value = "AEBDB54C343009E7B75F50D225DD0EF273B3709D4C7E6C6EE738C99AC592AA4E"
$wsA20.Range("A7").Value = "F6E368CBB00542D57BF5208510011056953F5BAB10DBF4F4AE07A315D7C7F3EC"
$wsA20.Range("B15").Select()

# --- "A19" becomes the active sheet/tab (was "A17"), with its selection on C33 ---
$wsA19.Activate()
$wsA19.Range("C33").Select()
